$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V content between the three reordered match-row pairs ---
$tmp = $ws.Range("F4:V4").Value2
$ws.Range("F4:V4").Value2 = $ws.Range("F5:V5").Value2
$ws.Range("F5:V5").Value2 = $tmp

$tmp = $ws.Range("F35:V35").Value2
$ws.Range("F35:V35").Value2 = $ws.Range("F36:V36").Value2
$ws.Range("F36:V36").Value2 = $tmp

$tmp = $ws.Range("F62:V62").Value2
$ws.Range("F62:V62").Value2 = $ws.Range("F63:V63").Value2
$ws.Range("F63:V63").Value2 = $tmp

# --- Append 4 new match rows (65-68), copying row-64 formatting first ---
$ws.Range("A64:V64").Copy()
$ws.Range("A65:V68").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 65
$ws.Range("A65").Value2 = 64
$ws.Range("B65").Value2 = "iran"
$ws.Range("C65").Value2 = "persian-gulf-pro-league"
$ws.Range("D65").Value2 = "2023-2024"
$ws.Range("E65").Value2 = 45240.52083333334
$ws.Range("F65").Value2 = "Gol Gohar"
$ws.Range("G65").Value2 = 2
$ws.Range("H65").Value2 = "Mes Rafsanjan"
$ws.Range("I65").Value2 = 0
$ws.Range("J65").Value2 = 2.54
$ws.Range("K65").Value2 = "09/11/2023 00:42"
$ws.Range("L65").Value2 = 2.89
$ws.Range("M65").Value2 = "10/11/2023 12:26"
$ws.Range("N65").Value2 = 2.66
$ws.Range("O65").Value2 = "09/11/2023 00:42"
$ws.Range("P65").Value2 = 2.45
$ws.Range("Q65").Value2 = "10/11/2023 12:26"
$ws.Range("R65").Value2 = 2.9
$ws.Range("S65").Value2 = "09/11/2023 00:42"
$ws.Range("T65").Value2 = 3.16
$ws.Range("U65").Value2 = "10/11/2023 12:26"
$ws.Range("V65").Value2 = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/gol-gohar-mes-rafsanjan/fJHJynpG/"

# Row 66
$ws.Range("A66").Value2 = 65
$ws.Range("B66").Value2 = "iran"
$ws.Range("C66").Value2 = "persian-gulf-pro-league"
$ws.Range("D66").Value2 = "2023-2024"
$ws.Range("E66").Value2 = 45240.52083333334
$ws.Range("F66").Value2 = "Malavan"
$ws.Range("G66").Value2 = 1
$ws.Range("H66").Value2 = "Esteghlal Khuzestan"
$ws.Range("I66").Value2 = 1
$ws.Range("J66").Value2 = 1.81
$ws.Range("K66").Value2 = "09/11/2023 00:42"
$ws.Range("L66").Value2 = 1.39
$ws.Range("M66").Value2 = "10/11/2023 12:29"
$ws.Range("N66").Value2 = 2.92
$ws.Range("O66").Value2 = "09/11/2023 00:42"
$ws.Range("P66").Value2 = 3.45
$ws.Range("Q66").Value2 = "10/11/2023 12:29"
$ws.Range("R66").Value2 = 4.57
$ws.Range("S66").Value2 = "09/11/2023 00:42"
$ws.Range("T66").Value2 = 7.84
$ws.Range("U66").Value2 = "10/11/2023 12:29"
$ws.Range("V66").Value2 = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/malavan-esteghlal-khuzestan/CfSExSU9/"

# Row 67
$ws.Range("A67").Value2 = 66
$ws.Range("B67").Value2 = "iran"
$ws.Range("C67").Value2 = "persian-gulf-pro-league"
$ws.Range("D67").Value2 = "2023-2024"
$ws.Range("E67").Value2 = 45240.52083333334
$ws.Range("F67").Value2 = "Paykan"
$ws.Range("G67").Value2 = 1
$ws.Range("H67").Value2 = "Shams Azar Qazvin"
$ws.Range("I67").Value2 = 1
$ws.Range("J67").Value2 = 2.68
$ws.Range("K67").Value2 = "09/11/2023 00:42"
$ws.Range("L67").Value2 = 4.59
$ws.Range("M67").Value2 = "10/11/2023 12:25"
$ws.Range("N67").Value2 = 2.64
$ws.Range("O67").Value2 = "09/11/2023 00:42"
$ws.Range("P67").Value2 = 2.73
$ws.Range("Q67").Value2 = "10/11/2023 12:25"
$ws.Range("R67").Value2 = 2.77
$ws.Range("S67").Value2 = "09/11/2023 00:42"
$ws.Range("T67").Value2 = 2.05
$ws.Range("U67").Value2 = "10/11/2023 12:25"
$ws.Range("V67").Value2 = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/paykan-shams-azar-qazvin/biOAw8F3/"

# Row 68
$ws.Range("A68").Value2 = 67
$ws.Range("B68").Value2 = "iran"
$ws.Range("C68").Value2 = "persian-gulf-pro-league"
$ws.Range("D68").Value2 = "2023-2024"
$ws.Range("E68").Value2 = 45240.59375
$ws.Range("F68").Value2 = "Foolad"
$ws.Range("G68").Value2 = 1
$ws.Range("H68").Value2 = "Aluminium Arak"
$ws.Range("I68").Value2 = 0
$ws.Range("J68").Value2 = 2.2
$ws.Range("K68").Value2 = "09/11/2023 02:42"
$ws.Range("L68").Value2 = 2.51
$ws.Range("M68").Value2 = "10/11/2023 14:14"
$ws.Range("N68").Value2 = 2.58
$ws.Range("O68").Value2 = "09/11/2023 02:42"
$ws.Range("P68").Value2 = 2.45
$ws.Range("Q68").Value2 = "10/11/2023 14:14"
$ws.Range("R68").Value2 = 3.68
$ws.Range("S68").Value2 = "09/11/2023 02:42"
$ws.Range("T68").Value2 = 3.78
$ws.Range("U68").Value2 = "10/11/2023 14:14"
$ws.Range("V68").Value2 = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/foolad-aluminium-arak/Y9GNz6aM/"

